$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Save clean style templates (outside the 1:24 working range) to use as a
# paste-formats source, since this runtime mis-resolves overlapping <col>
# style ranges for column B when a brand-new cell is created there.
$ws.Range("A3").Copy()
$ws.Range("A200").PasteSpecial(-4122)
$ws.Range("B3").Copy()
$ws.Range("B200").PasteSpecial(-4122)
$ws.Range("C3").Copy()
$ws.Range("C200").PasteSpecial(-4122)

# Wipe the existing table completely; it will be rebuilt from scratch below.
$ws.Rows("1:24").Delete()

# Template row 200 shifted up by 24 rows (the deleted range) -> row 176.
$tmplRow = 176

# Row 1
$ws.Range("B1").Value = "Ementa atual:"
$ws.Range("B" + $tmplRow).Copy()
$ws.Range("B1").PasteSpecial(-4122)
$ws.Range("C1").Value = "Ementa modificada (dados modificados em vermelho):"
$ws.Range("C" + $tmplRow).Copy()
$ws.Range("C1").PasteSpecial(-4122)

# Row 2
$ws.Range("B2").Value = "LOM3211"
$ws.Range("B" + $tmplRow).Copy()
$ws.Range("B2").PasteSpecial(-4122)
$ws.Range("C2").Value = "LOM3211"
$ws.Range("C" + $tmplRow).Copy()
$ws.Range("C2").PasteSpecial(-4122)

# Row 3
$ws.Range("A3").Value = "Nome:"
$ws.Range("A" + $tmplRow).Copy()
$ws.Range("A3").PasteSpecial(-4122)
$ws.Range("B3").Value = " Estruturas e Propriedades de Materiais"
$ws.Range("B" + $tmplRow).Copy()
$ws.Range("B3").PasteSpecial(-4122)
$ws.Range("C3").Value = " Estruturas e Propriedades de Materiais"
$ws.Range("C" + $tmplRow).Copy()
$ws.Range("C3").PasteSpecial(-4122)

# Row 4
$ws.Range("A4").Value = "Name:"
$ws.Range("A" + $tmplRow).Copy()
$ws.Range("A4").PasteSpecial(-4122)
$ws.Range("B4").Value = "Structures and Properties of Materials"
$ws.Range("B" + $tmplRow).Copy()
$ws.Range("B4").PasteSpecial(-4122)
$ws.Range("C4").Value = "Structures and Properties of Materials"
$ws.Range("C" + $tmplRow).Copy()
$ws.Range("C4").PasteSpecial(-4122)

# Row 5
$ws.Range("A5").Value = "Créditos-aula:"
$ws.Range("A" + $tmplRow).Copy()
$ws.Range("A5").PasteSpecial(-4122)
$ws.Range("B5").Value = "4"
$ws.Range("B" + $tmplRow).Copy()
$ws.Range("B5").PasteSpecial(-4122)
$ws.Range("C5").Value = "4"
$ws.Range("C" + $tmplRow).Copy()
$ws.Range("C5").PasteSpecial(-4122)

# Row 6
$ws.Range("A6").Value = "Créditos-trabalho"
$ws.Range("A" + $tmplRow).Copy()
$ws.Range("A6").PasteSpecial(-4122)
$ws.Range("B6").Value = "0"
$ws.Range("B" + $tmplRow).Copy()
$ws.Range("B6").PasteSpecial(-4122)
$ws.Range("C6").Value = "0"
$ws.Range("C" + $tmplRow).Copy()
$ws.Range("C6").PasteSpecial(-4122)

# Row 7
$ws.Range("A7").Value = "Carga horária:"
$ws.Range("A" + $tmplRow).Copy()
$ws.Range("A7").PasteSpecial(-4122)
$ws.Range("B7").Value = "60 h"
$ws.Range("B" + $tmplRow).Copy()
$ws.Range("B7").PasteSpecial(-4122)
$ws.Range("C7").Value = "60 h"
$ws.Range("C" + $tmplRow).Copy()
$ws.Range("C7").PasteSpecial(-4122)

# Row 8
$ws.Range("A8").Value = "Ativação:"
$ws.Range("A" + $tmplRow).Copy()
$ws.Range("A8").PasteSpecial(-4122)
$ws.Range("B8").Value = "01/01/2012"
$ws.Range("B" + $tmplRow).Copy()
$ws.Range("B8").PasteSpecial(-4122)
$ws.Range("C8").Value = "01/01/2012"
$ws.Range("C" + $tmplRow).Copy()
$ws.Range("C8").PasteSpecial(-4122)

# Row 9
$ws.Range("A9").Value = "Semestre ideal:"
$ws.Range("A" + $tmplRow).Copy()
$ws.Range("A9").PasteSpecial(-4122)
$ws.Range("B9").Value = "EF-7"
$ws.Range("B" + $tmplRow).Copy()
$ws.Range("B9").PasteSpecial(-4122)
$ws.Range("C9").Value = "EF-7"
$ws.Range("C" + $tmplRow).Copy()
$ws.Range("C9").PasteSpecial(-4122)

# Row 10
$ws.Range("A10").Value = "Objetivos:"
$ws.Range("A" + $tmplRow).Copy()
$ws.Range("A10").PasteSpecial(-4122)
$ws.Range("B10").Value = "984972 - Hugo Ricardo Zschommler Sandim"
$ws.Range("B" + $tmplRow).Copy()
$ws.Range("B10").PasteSpecial(-4122)
$ws.Range("C10").Value = "984972 - Hugo Ricardo Zschommler Sandim"
$ws.Range("C" + $tmplRow).Copy()
$ws.Range("C10").PasteSpecial(-4122)
$ws.Rows("10").RowHeight = 60

# Row 11
$ws.Range("A11").Value = "Objectives:"
$ws.Range("A" + $tmplRow).Copy()
$ws.Range("A11").PasteSpecial(-4122)
$ws.Rows("11").RowHeight = 60

# Row 12
$ws.Range("A12").Value = "Docentes responsáveis:"
$ws.Range("A" + $tmplRow).Copy()
$ws.Range("A12").PasteSpecial(-4122)

# Row 13
$ws.Range("A13").Value = "Programa resumido:"
$ws.Range("A" + $tmplRow).Copy()
$ws.Range("A13").PasteSpecial(-4122)
$ws.Range("B13").Value = "Semestral"
$ws.Range("B" + $tmplRow).Copy()
$ws.Range("B13").PasteSpecial(-4122)
$ws.Range("C13").Value = "Semestral"
$ws.Range("C" + $tmplRow).Copy()
$ws.Range("C13").PasteSpecial(-4122)
$ws.Rows("13").RowHeight = 60

# Row 14
$ws.Range("A14").Value = "Short syllabus:"
$ws.Range("A" + $tmplRow).Copy()
$ws.Range("A14").PasteSpecial(-4122)
$ws.Rows("14").RowHeight = 60

# Row 15
$ws.Range("A15").Value = "Programa:"
$ws.Range("A" + $tmplRow).Copy()
$ws.Range("A15").PasteSpecial(-4122)
$ws.Range("B15").Value = "01/01/2012"
$ws.Range("B" + $tmplRow).Copy()
$ws.Range("B15").PasteSpecial(-4122)
$ws.Range("C15").Value = "01/01/2012"
$ws.Range("C" + $tmplRow).Copy()
$ws.Range("C15").PasteSpecial(-4122)
$ws.Rows("15").RowHeight = 120

# Row 16
$ws.Range("A16").Value = "Syllabus:"
$ws.Range("A" + $tmplRow).Copy()
$ws.Range("A16").PasteSpecial(-4122)
$ws.Rows("16").RowHeight = 120

# Row 17
$ws.Range("A17").Value = "Avaliação:"
$ws.Range("A" + $tmplRow).Copy()
$ws.Range("A17").PasteSpecial(-4122)

# Row 18
$ws.Range("A18").Value = "Método:"
$ws.Range("A" + $tmplRow).Copy()
$ws.Range("A18").PasteSpecial(-4122)
$ws.Range("B18").Value = "984972 - Hugo Ricardo Zschommler Sandim"
$ws.Range("B" + $tmplRow).Copy()
$ws.Range("B18").PasteSpecial(-4122)
$ws.Range("C18").Value = "984972 - Hugo Ricardo Zschommler Sandim"
$ws.Range("C" + $tmplRow).Copy()
$ws.Range("C18").PasteSpecial(-4122)
$ws.Rows("18").RowHeight = 60

# Row 19
$ws.Range("A19").Value = "Critério:"
$ws.Range("A" + $tmplRow).Copy()
$ws.Range("A19").PasteSpecial(-4122)
$ws.Range("B19").Value = "Aulas expositivas, seminários e exercícios comentados."
$ws.Range("B" + $tmplRow).Copy()
$ws.Range("B19").PasteSpecial(-4122)
$ws.Range("C19").Value = "Aulas expositivas, seminários e exercícios comentados."
$ws.Range("C" + $tmplRow).Copy()
$ws.Range("C19").PasteSpecial(-4122)
$ws.Rows("19").RowHeight = 60

# Row 20
$ws.Range("A20").Value = "Norma de recuperação:"
$ws.Range("A" + $tmplRow).Copy()
$ws.Range("A20").PasteSpecial(-4122)
$ws.Range("B20").Value = "Média aritmética de duas provas sendo a primeira com peso 1 e a segunda com peso 2."
$ws.Range("B" + $tmplRow).Copy()
$ws.Range("B20").PasteSpecial(-4122)
$ws.Range("C20").Value = "Média aritmética de duas provas sendo a primeira com peso 1 e a segunda com peso 2."
$ws.Range("C" + $tmplRow).Copy()
$ws.Range("C20").PasteSpecial(-4122)
$ws.Rows("20").RowHeight = 60

# Row 21
$ws.Range("A21").Value = "Bibliografia:"
$ws.Range("A" + $tmplRow).Copy()
$ws.Range("A21").PasteSpecial(-4122)
$ws.Range("B21").Value = "Aplicação de uma prova escrita dentro do prazo regimental antes do início do próximo semestre letivo. A nota da segunda avaliação será a média aritmética entre a nota da prova de recuperação e a nota final da primeira avaliação"
$ws.Range("B" + $tmplRow).Copy()
$ws.Range("B21").PasteSpecial(-4122)
$ws.Range("C21").Value = "Aplicação de uma prova escrita dentro do prazo regimental antes do início do próximo semestre letivo. A nota da segunda avaliação será a média aritmética entre a nota da prova de recuperação e a nota final da primeira avaliação"
$ws.Range("C" + $tmplRow).Copy()
$ws.Range("C21").PasteSpecial(-4122)
$ws.Rows("21").RowHeight = 120

# Row 22
$ws.Range("A22").Value = "Requisitos:"
$ws.Range("A" + $tmplRow).Copy()
$ws.Range("A22").PasteSpecial(-4122)

# Row 23
$ws.Range("B23").Value = "LOM3016 -  Introdução à  Ciência dos Materiais  (Requisito)`n"
$ws.Range("B" + $tmplRow).Copy()
$ws.Range("B23").PasteSpecial(-4122)
$ws.Range("C23").Value = "LOM3016 -  Introdução à  Ciência dos Materiais  (Requisito)`n"
$ws.Range("C" + $tmplRow).Copy()
$ws.Range("C23").PasteSpecial(-4122)
$ws.Rows("23").RowHeight = 30

# Remove the temporary style-template row.
$ws.Rows($tmplRow).Delete()

$ws.Range("A1").Select()
Write-Output "done"
